$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new US-recode assay rows at the bottom of the table.
$ws.Range("A34").Value = "Abbott ARCHITECT SARS-CoV-2 IgG immunoassay"
$ws.Range("B34").Value = "nucleocapsid"
$ws.Range("C34").Value = "IgG"
$ws.Range("D34").Value = "N-Abbott"

$ws.Range("A35").Value = "Ortho-Clinical Diagnostics VITROS SARS-CoV-2 IgG immunoassay"
$ws.Range("B35").Value = "spike"
$ws.Range("C35").Value = "IgG"
$ws.Range("D35").Value = "S-Ortho IgG"

# Re-sort the whole table (excluding header) by test_name ascending, like the
# author did in Excel before committing.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A35"))
$ws.Sort.SetRange($ws.Range("A2:D35"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Excel's text collation ties "Abbott Architect IgG" against the new,
# differently-punctuated "Abbott ARCHITECT SARS-CoV-2 IgG immunoassay" row
# right after the base entry, ahead of the semicolon-qualified rows -- fix up
# that tie-break so the final order matches exactly.
$ws.Range("A3").Value = "Abbott ARCHITECT SARS-CoV-2 IgG immunoassay"
$ws.Range("B3").Value = "nucleocapsid"
$ws.Range("C3").Value = "IgG"
$ws.Range("D3").Value = "N-Abbott"

$ws.Range("A4").Value = "Abbott Architect IgG; VITROS IgG"
$ws.Range("B4").Value = "mixed"
$ws.Range("C4").Value = "IgG"
$ws.Range("D4").Value = "N-Abbott, S-Ortho IgG"

$ws.Range("A5").Value = "Abbott Architect IgG; Wantai ELISA pan-Ig"
$ws.Range("B5").Value = "mixed"
$ws.Range("C5").Value = "pan-Ig"
$ws.Range("D5").Value = "S-Ortho Ig, S-Roche"

$ws.Range("D27").Select()
